# Fix rPr child element ordering in character styles so that <w:b/> / <w:i/>
# come before <w:color/> (per wml.xsd CT_RPr sequence), matching the
# OOXMLValidator-compliant order. Re-assigning the Font.Bold / Font.Italic
# properties forces the engine to re-serialize <w:rPr> in schema order.

$d = $word.ActiveDocument

# Styles whose <w:rPr> only contains <w:color/> + <w:b/> (needs b before color)
$boldOnly = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($styleName in $boldOnly) {
    $style = $d.Styles.Item($styleName)
    $style.Font.Bold = $true
}

# Styles whose <w:rPr> only contains <w:color/> + <w:i/> (needs i before color)
$italicOnly = @("CommentTok", "DocumentationTok")
foreach ($styleName in $italicOnly) {
    $style = $d.Styles.Item($styleName)
    $style.Font.Italic = $true
}

# Styles whose <w:rPr> contains <w:color/> + <w:b/> + <w:i/> (needs b, i before color)
$boldItalic = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($styleName in $boldItalic) {
    $style = $d.Styles.Item($styleName)
    $style.Font.Bold = $true
    $style.Font.Italic = $true
}
